$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'245.01"
$ws.Range("G2").Value = "'18"

# Row 3
$ws.Range("G3").Value = "'18"

# Row 4
$ws.Range("D4").Value = "'5.390"
$ws.Range("G4").Value = "'18"

# Row 5
$ws.Range("G5").Value = "'18"

# Row 6
$ws.Range("G6").Value = "'18"

# Row 7
$ws.Range("D7").Value = "'0.8152"
$ws.Range("G7").Value = "'18"

# Row 8
$ws.Range("D8").Value = "'0.9297"
$ws.Range("G8").Value = "'18"

# Row 9
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "'0.1435"
$ws.Range("E9").Value = "8WazirXWRX"
$ws.Range("G9").Value = "'18"

# Row 10
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").Value = "'0.07430"
$ws.Range("E10").Value = "9MandalaExchangeTokenMDX"
$ws.Range("G10").Value = "'18"

# Row 11
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "'0.03501"
$ws.Range("E11").Value = "10LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("G11").Value = "'18"

# Row 12
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03068"
$ws.Range("E12").Value = "11BitrueCoinBTR"
$ws.Range("G12").Value = "'18"

# Row 13
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09412"
$ws.Range("E13").Value = "12BitMartTokenBMX"
$ws.Range("G13").Value = "'18"

# Row 14
$ws.Range("B14").Value = "MCDex"
$ws.Range("C14").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D14").Value = "'4.006"
$ws.Range("E14").Value = "13MCDexMCB"
$ws.Range("G14").Value = "'18"

# Row 15
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001603"
$ws.Range("E15").Value = "14BitForexTokenBF"
$ws.Range("G15").Value = "'18"

# Row 16
$ws.Range("B16").Value = "CoinExToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D16").Value = "'0.04812"
$ws.Range("E16").Value = "15CoinExTokenCET"
$ws.Range("G16").Value = "'18"

# Row 17
$ws.Range("B17").Value = "One"
$ws.Range("C17").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D17").Value = "'0.0005941"
$ws.Range("E17").Value = "16OneONE"
$ws.Range("G17").Value = "'18"

# Row 18
$ws.Range("D18").Value = "'0.005404"
$ws.Range("G18").Value = "'18"

# Row 19
$ws.Range("D19").Value = "'0.004159"
$ws.Range("G19").Value = "'18"

# Row 20
$ws.Range("D20").Value = "'0.0009901"
$ws.Range("G20").Value = "'18"

# Row 21
$ws.Range("D21").Value = "'3.668"
$ws.Range("G21").Value = "'18"

# Row 22
$ws.Range("D22").Value = "'6.430"
$ws.Range("G22").Value = "'18"

# Row 23
$ws.Range("D23").Value = "'2.181"
$ws.Range("G23").Value = "'18"

# Row 24
$ws.Range("D24").Value = "'0.3252"
$ws.Range("G24").Value = "'18"

# Row 25
$ws.Range("G25").Value = "'18"

# Row 26
$ws.Range("D26").Value = "'0.00007000"
$ws.Range("G26").Value = "'18"

# Row 27
$ws.Range("D27").Value = "'0.0002900"
$ws.Range("G27").Value = "'18"

# Row 28
$ws.Range("G28").Value = "'18"

# Row 29
$ws.Range("G29").Value = "'18"

# Row 30
$ws.Range("G30").Value = "'18"

# Row 31
$ws.Range("G31").Value = "'18"

# Row 32
$ws.Range("G32").Value = "'18"

# Row 33
$ws.Range("G33").Value = "'18"

# Row 34
$ws.Range("G34").Value = "'18"

# Row 35
$ws.Range("G35").Value = "'18"

# Row 36
$ws.Range("G36").Value = "'18"

# Row 37
$ws.Range("G37").Value = "'18"

# Row 38
$ws.Range("G38").Value = "'18"

# Row 39
$ws.Range("G39").Value = "'18"

# Row 40
$ws.Range("D40").Value = "'0.04022"
$ws.Range("G40").Value = "'18"

# Row 41
$ws.Range("D41").Value = "'0.006370"
$ws.Range("G41").Value = "'18"

# Row 42
$ws.Range("G42").Value = "'18"

# Row 43
$ws.Range("D43").Value = "'0.002720"
$ws.Range("G43").Value = "'18"

# Row 44
$ws.Range("D44").Value = "'0.006603"
$ws.Range("E44").Value = "43LocalTradersLCTBestin24h"
$ws.Range("G44").Value = "'18"

# Row 45
$ws.Range("D45").Value = "'0.00005284"
$ws.Range("G45").Value = "'18"

# Row 46
$ws.Range("G46").Value = "'18"

# Row 47
$ws.Range("D47").Value = "'0.8601"
$ws.Range("G47").Value = "'18"

# Row 48
$ws.Range("D48").Value = "'0.002468"
$ws.Range("G48").Value = "'18"

# Row 49
$ws.Range("G49").Value = "'18"

# Row 50
$ws.Range("G50").Value = "'18"

# Row 51
$ws.Range("G51").Value = "'18"
